# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the existing "Late" / "Outstanding" data one column to
# the right (N->O, O->P, P->Q), and switch the active sheet/selection from
# "Transactions" back to "Repayment Schedule".

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column at N (column 14); everything from N onward shifts
# right by one column (N->O, O->P, P->Q).
$schedule.Columns.Item(14).Insert() | Out-Null

# Make "Repayment Schedule" the active sheet again (it was "Transactions"),
# and leave the selection where the author's last edit left it.
$schedule.Activate() | Out-Null
$schedule.Range("K25").Select() | Out-Null
